# IN-719 make all entities work with table defs
# Adds a new "table_definitions" sheet in front of the existing sheets and
# populates it with the table-definition metadata for crec_persons, and
# tweaks a couple of pre-existing layout/formatting details.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "table_definitions" sheet as the first sheet -------
$defSheet = $wb.Worksheets.Add()
$defSheet.Name = "table_definitions"
$defSheet.Move($wb.Worksheets.Item(1))

# Header row
$headers = @( `
    "mapping_file_name", `
    "entity_name", `
    "required_entities", `
    "destination_table_name", `
    "table_type", `
    "source_table_name", `
    "casrec_conditions", `
    "source_table_additional_columns" `
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $defSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row describing the crec_persons entity
$defSheet.Range("A2").Value = "crec_persons"
$defSheet.Range("B2").Value = "crec"
$defSheet.Range("C2").Value = "client"
$defSheet.Range("D2").Value = "persons"
$defSheet.Range("E2").Value = "data"
$defSheet.Range("F2").Value = "crec"

$defSheet.Rows("1:2").Select() | Out-Null

# --- 2. Tighten up the crec_lookup sheet's middle rows (2-5) --------------
$lookup = $wb.Worksheets.Item("crec_lookup")
$lookup.Rows("2:5").RowHeight = 15

$defSheet.Activate() | Out-Null
